$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = "'0"
$ws.Range("N15").Value = "-40"
$ws.Range("C16").Value = "1"
$ws.Range("D16").Value = "1"
$ws.Range("E16").Value = "0"
$ws.Range("F16").Value = "4"
$ws.Range("G16").Value = "10"
$ws.Range("H16").Value = "-60"
$ws.Range("I16").Value = "22"
$ws.Range("J16").Value = "32"
$ws.Range("K16").Value = "-31.25"
$ws.Range("L16").Value = "-4.347826086956"
$ws.Range("M16").Value = "-51.111111111111"
$ws.Range("N16").Value = "-84.397163120567"
$ws.Range("C17").Value = "1"
$ws.Range("D17").Value = "3"
$ws.Range("E17").Value = "-66.666666666666"
$ws.Range("F17").Value = "9"
$ws.Range("H17").Value = "-10"
$ws.Range("I17").Value = "23"
$ws.Range("J17").Value = "27"
$ws.Range("K17").Value = "-14.814814814814"
$ws.Range("L17").Value = "-36.111111111111"
$ws.Range("M17").Value = "4.545454545454"
$ws.Range("N17").Value = "-68.918918918918"
$ws.Range("C18").Value = "4"
$ws.Range("D18").Value = "11"
$ws.Range("E18").Value = "-63.636363636363"
$ws.Range("F18").Value = "12"
$ws.Range("G18").Value = "24"
$ws.Range("H18").Value = "-50"
$ws.Range("I18").Value = "24"
$ws.Range("J18").Value = "53"
$ws.Range("K18").Value = "-54.716981132075"
$ws.Range("L18").Value = "-27.272727272727"
$ws.Range("M18").Value = "-7.692307692307"
$ws.Range("N18").Value = "-89.473684210526"
$ws.Range("C19").Value = "10"
$ws.Range("D19").Value = "7"
$ws.Range("E19").Value = "42.857142857142"
$ws.Range("F19").Value = "33"
$ws.Range("G19").Value = "30"
$ws.Range("H19").Value = "10"
$ws.Range("I19").Value = "102"
$ws.Range("J19").Value = "92"
$ws.Range("K19").Value = "10.869565217391"
$ws.Range("L19").Value = "78.947368421052"
$ws.Range("M19").Value = "10.869565217391"
$ws.Range("N19").Value = "-47.150259067357"
$ws.Range("C20").Value = "2"
$ws.Range("D20").Value = "1"
$ws.Range("E20").Value = "100"
$ws.Range("F20").Value = "5"
$ws.Range("G20").Value = "6"
$ws.Range("H20").Value = "-16.666666666666"
$ws.Range("I20").Value = "22"
$ws.Range("J20").Value = "10"
$ws.Range("K20").Value = "120"
$ws.Range("L20").Value = "144.444444444444"
$ws.Range("M20").Value = "175"
$ws.Range("N20").Value = "-87.134502923976"
$ws.Range("C21").Value = "18"
$ws.Range("D21").Value = "23"
$ws.Range("E21").Value = "-21.739130434782"
$ws.Range("F21").Value = "65"
$ws.Range("G21").Value = "80"
$ws.Range("H21").Value = "-18.75"
$ws.Range("I21").Value = "198"
$ws.Range("J21").Value = "215"
$ws.Range("K21").Value = "-7.906976744186"
$ws.Range("L21").Value = "24.528301886792"
$ws.Range("M21").Value = "1.538461538461"
$ws.Range("N21").Value = "-75.764993880049"
$ws.Range("D22").Value = "1"
$ws.Range("E22").Value = "-100"
$ws.Range("J22").Value = "8"
$ws.Range("K22").Value = "-62.5"
$ws.Range("L22").Value = "-25"
$ws.Range("C23").Value = "'0"
$ws.Range("D23").Value = "2"
$ws.Range("E23").Value = "-100"
$ws.Range("F23").Value = "5"
$ws.Range("G23").Value = "9"
$ws.Range("H23").Value = "-44.444444444444"
$ws.Range("J23").Value = "20"
$ws.Range("K23").Value = "-10"
$ws.Range("L23").Value = "-45.454545454545"
$ws.Range("M23").Value = "-14.285714285714"
$ws.Range("C24").Value = "31"
$ws.Range("D24").Value = "47"
$ws.Range("E24").Value = "-34.042553191489"
$ws.Range("F24").Value = "117"
$ws.Range("G24").Value = "153"
$ws.Range("H24").Value = "-23.529411764705"
$ws.Range("I24").Value = "387"
$ws.Range("J24").Value = "362"
$ws.Range("K24").Value = "6.906077348066"
$ws.Range("L24").Value = "118.64406779661"
$ws.Range("M24").Value = "106.951871657754"
$ws.Range("C25").Value = "6"
$ws.Range("D25").Value = "2"
$ws.Range("E25").Value = "200"
$ws.Range("F25").Value = "18"
$ws.Range("G25").Value = "16"
$ws.Range("H25").Value = "12.5"
$ws.Range("I25").Value = "50"
$ws.Range("J25").Value = "48"
$ws.Range("K25").Value = "4.166666666666"
$ws.Range("L25").Value = "42.857142857142"
$ws.Range("M25").Value = "-21.875"
$ws.Range("G26").Value = "'0"
$ws.Range("H26").Value = "'***.*"
$ws.Range("L26").Value = "33.333333333333"
$ws.Range("C27").Value = "4"
$ws.Range("F27").Value = "6"
$ws.Range("G27").Value = "3"
$ws.Range("H27").Value = "100"
$ws.Range("I27").Value = "11"
$ws.Range("K27").Value = "37.5"
$ws.Range("L27").Value = "10"
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "'***.*"
$ws.Range("G30").Value = "5"

# Fix up styles for cells that changed numeric/text type (one at a time to avoid multi-area paste issues)
$ws.Range("A14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("M14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("M14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
